# "Generate Report for Handoff"
# The localization status report is regenerated: the Overview/zh-cn/de-de
# status moves from "In Translation" to "Ready for handoff", the associated
# timestamps advance a little, and the "latest handback" URL picks up a new
# commit SHA. The Status / zh-cn / de-de columns also grow a bit wider to
# fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2017-02-09 14:18:34"

# --- zh-cn sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime, R2 = Error Detail
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2017-02-09 14:18:17"
$wsZhCn.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0729dc005dfb2c635e2cf1a74b23e5cacd7ace06/e2e/ae7f0526-159b-4eaf-aafd-6e77a2be2935.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2ecb50fad0d4fef7d9606c8ba7986f5bae1df96/e2e/ae7f0526-159b-4eaf-aafd-6e77a2be2935.md."

# --- de-de sheet --------------------------------------------------------
# C2 = Status, H2 = Latest Handback DateTime, R2 = Error Detail
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2017-02-09 14:18:34"
$wsDeDe.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0729dc005dfb2c635e2cf1a74b23e5cacd7ace06/e2e/ae7f0526-159b-4eaf-aafd-6e77a2be2935.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2ecb50fad0d4fef7d9606c8ba7986f5bae1df96/e2e/ae7f0526-159b-4eaf-aafd-6e77a2be2935.md."

# --- Column widths --------------------------------------------------
# The longer "Ready for handoff" text makes the Status / zh-cn / de-de
# columns autofit a bit wider (target stored OOXML width ~17.216 chars).
$newWidth = 17.2159881591797 - 0.8333333333333334
$wsOverview.Range("E:E").ColumnWidth = $newWidth
$wsOverview.Range("F:F").ColumnWidth = $newWidth
$wsZhCn.Range("C:C").ColumnWidth = $newWidth
$wsDeDe.Range("C:C").ColumnWidth = $newWidth
